# Add the "Visual Studio ALM Rangers" attribution textbox to slide 1.
#
# OOXML target (from diff): a new p:sp placed after the last p:pic in the
# slide's spTree, positioned at off (508000, 8661400) EMU with extent
# (5835650, 338554) EMU, noFill shape background, word-wrapped/auto-fit
# body, and three text runs (8pt) forming one sentence.
#
# PowerPoint's Shapes.AddTextbox (and the xfrm it writes) works in points,
# not EMU, so convert: 1 pt = 12700 EMU.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$left   = 508000  / 12700
$top    = 8661400 / 12700
$width  = 5835650 / 12700
$height = 338554  / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)

$tr = $tb.TextFrame.TextRange
$tr.Text = "This content was created by the Visual Studio ALM Rangers, a special group with members from the Visual Studio Product Team"
$tr.Font.Size = 8

$r2 = $tr.InsertAfter(", Microsoft ")
$r2.Font.Size = 8

$r3 = $tr.InsertAfter("Services, Microsoft Most Valuable Professionals (MVPs) and Visual Studio Community Leads.")
$r3.Font.Size = 8

$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1

$tb.Fill.Visible = $false
